$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35, pushing existing rows 35-66 down to 36-67.
$ws.Rows("35:35").Insert()

# Copy date cell style (s="2") from the row above (D34) onto the new D35 cell.
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 35 with the new weekly record.
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 45033
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112001
$ws.Cells.Item(35, 7).Value = "Berenjena"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 50
$ws.Cells.Item(35, 11).Value = 12000
$ws.Cells.Item(35, 12).Value = 12000
$ws.Cells.Item(35, 13).Value = 12000
$ws.Cells.Item(35, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 200
$ws.Cells.Item(35, 17).Value = 60
$ws.Cells.Item(35, 18).Value = "Hortaliza"
